$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper note: D-column cell values that look like plain decimal numbers
# (single dot) are entered with a leading apostrophe so Excel stores them
# as text, matching the source workbook's text-typed Price column. Values
# with multiple dots (e.g. "66.384.54") or other non-numeric characters
# are already safe to assign directly.

$ws.Range("D2").Value = "66.384.54"
$ws.Range("E2").Value = "  +1.43%  "

$ws.Range("D3").Value = "3.687.45"
$ws.Range("E3").Value = "  +4.45%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'418.53"
$ws.Range("E5").Value = "  -0.77%  "

$ws.Range("D6").Value = "'130.09"
$ws.Range("E6").Value = "  -3.52%  "

$ws.Range("D7").Value = "3.679.07"
$ws.Range("E7").Value = "  +4.45%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D10").Value = "'0.777"
$ws.Range("E10").Value = "  -1.93%  "

$ws.Range("D11").Value = "'0.178"
$ws.Range("E11").Value = "  +7.29%  "

$ws.Range("D12").Value = "'0.0000393"
$ws.Range("E12").Value = "  +45.34%  "

$ws.Range("D13").Value = "'43.07"
$ws.Range("E13").Value = "  -1.59%  "

$ws.Range("E14").Value = "  +5.19%  "

$ws.Range("D15").Value = "4.268.83"
$ws.Range("E15").Value = "  +4.34%  "

$ws.Range("E16").Value = "  -0.94%  "

$ws.Range("D17").Value = "'20.57"
$ws.Range("E17").Value = "  -0.75%  "

$ws.Range("D18").Value = "3.678.10"
$ws.Range("E18").Value = "  +3.67%  "

$ws.Range("D19").Value = "'13.31"
$ws.Range("E19").Value = "  +5.37%  "

$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("D21").Value = "66.368.02"
$ws.Range("E21").Value = "  +1.61%  "

$ws.Range("D22").Value = "'444.07"
$ws.Range("E22").Value = "  -3.55%  "

$ws.Range("D23").Value = "'16.56"
$ws.Range("E23").Value = "  +23.40%  "

$ws.Range("D24").Value = "'89.87"
$ws.Range("E24").Value = "  -2.36%  "

$ws.Range("D25").Value = "'3.14"
$ws.Range("E25").Value = "  -3.01%  "

$ws.Range("D26").Value = "'37.16"
$ws.Range("E26").Value = "  +7.65%  "

$ws.Range("E27").Value = "  +0.65%  "

$ws.Range("E28").Value = "  -1.33%  "

$ws.Range("D29").Value = "'5.02"
$ws.Range("E29").Value = "  +4.06%  "

$ws.Range("E30").Value = "  +8.45%  "

$ws.Range("D31").Value = "'12.71"
$ws.Range("E31").Value = "  +0.99%  "

$ws.Range("D32").Value = "'2.78"
$ws.Range("E32").Value = "  -1.90%  "

$ws.Range("D33").Value = "'7.30"
$ws.Range("E33").Value = "  -3.75%  "

$ws.Range("D34").Value = "'0.165"
$ws.Range("E34").Value = "  +1.38%  "

$ws.Range("D35").Value = "'41.64"
$ws.Range("E35").Value = "  +3.09%  "

$ws.Range("D36").Value = "'57.31"
$ws.Range("E36").Value = "  -0.70%  "

$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").Value = "'0.0492"
$ws.Range("E38").Value = "  -2.94%  "

$ws.Range("E39").Value = "  +33.29%  "

$ws.Range("D40").Value = "0.0₃0725"
$ws.Range("E40").Value = "  -3.65%  "

$ws.Range("D41").Value = "'0.149"
$ws.Range("E41").Value = "  +2.78%  "

$ws.Range("D42").Value = "'29.52"
$ws.Range("E42").Value = "  +34.31%  "

$ws.Range("D43").Value = "'0.996"
$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("E44").Value = "  +2.27%  "

$ws.Range("D45").Value = "'148.71"
$ws.Range("E45").Value = "  +1.60%  "

$ws.Range("E46").Value = "  +4.47%  "

# Row 47: was WEMIXToken, now NEARProtocol
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'4.38"
$ws.Range("E47").Value = "  -4.39%  "

# Row 48: was NEARProtocol, now WEMIXToken
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'2.66"
$ws.Range("E48").Value = "  -5.12%  "

$ws.Range("D49").Value = "'2.89"
$ws.Range("E49").Value = "  -7.56%  "

$ws.Range("D50").Value = "'0.306"
$ws.Range("E50").Value = "  -3.23%  "

# Row 51: was Cronos, now ApeXProtocol
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "'3.02"
$ws.Range("E51").Value = "  +22.08%  "
